# This script updates odds data cells in the "Jogos da Semana" worksheet
# to reflect refreshed values scraped from FlashScore, per the commit
# "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.88
$ws.Cells.Item(2, 20).Value = 5.5
$ws.Cells.Item(2, 23).Value = 34
$ws.Cells.Item(2, 24).Value = 41
$ws.Cells.Item(2, 31).Value = 12
$ws.Cells.Item(2, 35).Value = 51

# Row 4
$ws.Cells.Item(4, 7).Value = 2.6
$ws.Cells.Item(4, 9).Value = 2.4
$ws.Cells.Item(4, 16).Value = 1.32
$ws.Cells.Item(4, 17).Value = 3.2
$ws.Cells.Item(4, 20).Value = 10.25
$ws.Cells.Item(4, 21).Value = 15.5
$ws.Cells.Item(4, 22).Value = 10.5
$ws.Cells.Item(4, 23).Value = 30
$ws.Cells.Item(4, 24).Value = 21
$ws.Cells.Item(4, 31).Value = 14
$ws.Cells.Item(4, 32).Value = 10
$ws.Cells.Item(4, 33).Value = 27
$ws.Cells.Item(4, 34).Value = 19

# Row 5
$ws.Cells.Item(5, 8).Value = 3.1
$ws.Cells.Item(5, 9).Value = 2.5
$ws.Cells.Item(5, 11).Value = 6.7
$ws.Cells.Item(5, 12).Value = 1.38
$ws.Cells.Item(5, 13).Value = 2.87
$ws.Cells.Item(5, 14).Value = 2.15
$ws.Cells.Item(5, 15).Value = 1.65
$ws.Cells.Item(5, 16).Value = 1.45
$ws.Cells.Item(5, 17).Value = 2.62
$ws.Cells.Item(5, 18).Value = 1.82
$ws.Cells.Item(5, 19).Value = 1.88
$ws.Cells.Item(5, 20).Value = 7.8
$ws.Cells.Item(5, 24).Value = 30
$ws.Cells.Item(5, 26).Value = 6.7
$ws.Cells.Item(5, 28).Value = 15.5
$ws.Cells.Item(5, 29).Value = 80
$ws.Cells.Item(5, 30).Value = 7.5
$ws.Cells.Item(5, 31).Value = 13
$ws.Cells.Item(5, 32).Value = 10
$ws.Cells.Item(5, 34).Value = 23
$ws.Cells.Item(5, 35).Value = 35

# Row 9
$ws.Cells.Item(9, 34).Value = 32

# Row 10
$ws.Cells.Item(10, 7).Value = 2.1
$ws.Cells.Item(10, 8).Value = 3.25
$ws.Cells.Item(10, 9).Value = 3.5
$ws.Cells.Item(10, 12).Value = 1.5
$ws.Cells.Item(10, 13).Value = 2.5
$ws.Cells.Item(10, 21).Value = 8.5
$ws.Cells.Item(10, 23).Value = 19
$ws.Cells.Item(10, 24).Value = 21
$ws.Cells.Item(10, 30).Value = 7.5
$ws.Cells.Item(10, 31).Value = 17
$ws.Cells.Item(10, 35).Value = 51

# Row 11
$ws.Cells.Item(11, 19).Value = 2.05
$ws.Cells.Item(11, 26).Value = 14
$ws.Cells.Item(11, 30).Value = 16.5
$ws.Cells.Item(11, 32).Value = 14.5
$ws.Cells.Item(11, 36).Value = 300

# Row 17
$ws.Cells.Item(17, 8).Value = 3.6
$ws.Cells.Item(17, 9).Value = 5.25
$ws.Cells.Item(17, 10).Value = 1.07
$ws.Cells.Item(17, 11).Value = 9
$ws.Cells.Item(17, 12).Value = 1.36
$ws.Cells.Item(17, 13).Value = 3
$ws.Cells.Item(17, 14).Value = 2.2
$ws.Cells.Item(17, 15).Value = 1.65
$ws.Cells.Item(17, 16).Value = 1.44
$ws.Cells.Item(17, 17).Value = 2.63
$ws.Cells.Item(17, 18).Value = 2.1
$ws.Cells.Item(17, 19).Value = 1.67
$ws.Cells.Item(17, 20).Value = 5.5
$ws.Cells.Item(17, 21).Value = 7
$ws.Cells.Item(17, 22).Value = 9
$ws.Cells.Item(17, 24).Value = 15
$ws.Cells.Item(17, 25).Value = 34
$ws.Cells.Item(17, 26).Value = 8.5
$ws.Cells.Item(17, 28).Value = 21
$ws.Cells.Item(17, 29).Value = 67
$ws.Cells.Item(17, 30).Value = 12
$ws.Cells.Item(17, 32).Value = 19
$ws.Cells.Item(17, 35).Value = 51

# Row 19
$ws.Cells.Item(19, 7).Value = 2.8
$ws.Cells.Item(19, 9).Value = 2.3
$ws.Cells.Item(19, 12).Value = 1.22
$ws.Cells.Item(19, 13).Value = 4
$ws.Cells.Item(19, 14).Value = 1.75
$ws.Cells.Item(19, 15).Value = 2.05
$ws.Cells.Item(19, 31).Value = 12
$ws.Cells.Item(19, 34).Value = 17

# Row 20
$ws.Cells.Item(20, 8).Value = 3.2
$ws.Cells.Item(20, 10).Value = 1.1
$ws.Cells.Item(20, 11).Value = 7
$ws.Cells.Item(20, 12).Value = 1.44
$ws.Cells.Item(20, 13).Value = 2.63
$ws.Cells.Item(20, 22).Value = 9.5
$ws.Cells.Item(20, 26).Value = 7
$ws.Cells.Item(20, 30).Value = 9

# Row 21
$ws.Cells.Item(21, 11).Value = 5.5
$ws.Cells.Item(21, 12).Value = 1.47
$ws.Cells.Item(21, 13).Value = 2.45
$ws.Cells.Item(21, 14).Value = 2.4
$ws.Cells.Item(21, 20).Value = 5.9
$ws.Cells.Item(21, 21).Value = 9.25
$ws.Cells.Item(21, 22).Value = 9
$ws.Cells.Item(21, 23).Value = 21
$ws.Cells.Item(21, 24).Value = 20
$ws.Cells.Item(21, 25).Value = 37
$ws.Cells.Item(21, 26).Value = 5.5
$ws.Cells.Item(21, 30).Value = 8
$ws.Cells.Item(21, 31).Value = 17.5

# Row 23
$ws.Cells.Item(23, 12).Value = 1.4
$ws.Cells.Item(23, 13).Value = 2.75
$ws.Cells.Item(23, 14).Value = 2.25
$ws.Cells.Item(23, 15).Value = 1.62

# Row 27
$ws.Cells.Item(27, 7).Value = 1.95
$ws.Cells.Item(27, 8).Value = 3.15
$ws.Cells.Item(27, 9).Value = 3.85
$ws.Cells.Item(27, 12).Value = 1.38
$ws.Cells.Item(27, 13).Value = 2.6
$ws.Cells.Item(27, 14).Value = 2.1
$ws.Cells.Item(27, 15).Value = 1.57
$ws.Cells.Item(27, 16).Value = 1.47
$ws.Cells.Item(27, 17).Value = 2.32
$ws.Cells.Item(27, 18).Value = 1.9
$ws.Cells.Item(27, 19).Value = 1.72
$ws.Cells.Item(27, 20).Value = 6.2
$ws.Cells.Item(27, 21).Value = 8.5
$ws.Cells.Item(27, 22).Value = 8.5
$ws.Cells.Item(27, 23).Value = 17
$ws.Cells.Item(27, 24).Value = 17
$ws.Cells.Item(27, 25).Value = 32
$ws.Cells.Item(27, 26).Value = 7.7
$ws.Cells.Item(27, 27).Value = 6.2
$ws.Cells.Item(27, 28).Value = 16.5
$ws.Cells.Item(27, 30).Value = 9
$ws.Cells.Item(27, 31).Value = 19.5
$ws.Cells.Item(27, 32).Value = 13.5
$ws.Cells.Item(27, 33).Value = 60
$ws.Cells.Item(27, 34).Value = 40
$ws.Cells.Item(27, 35).Value = 55
$ws.Cells.Item(27, 36).Value = 900

# Row 31
$ws.Cells.Item(31, 30).Value = 11
$ws.Cells.Item(31, 33).Value = 41

# Row 33
$ws.Cells.Item(33, 10).Value = 1.1
$ws.Cells.Item(33, 18).Value = 2.15
$ws.Cells.Item(33, 19).Value = 1.62
$ws.Cells.Item(33, 20).Value = 5.5
$ws.Cells.Item(33, 33).Value = 55

# Row 34
$ws.Cells.Item(34, 7).Value = 1.39
$ws.Cells.Item(34, 8).Value = 4.4
$ws.Cells.Item(34, 9).Value = 6.6
$ws.Cells.Item(34, 18).Value = 1.98
$ws.Cells.Item(34, 19).Value = 1.75
$ws.Cells.Item(34, 20).Value = 6.9
$ws.Cells.Item(34, 21).Value = 6.5
$ws.Cells.Item(34, 22).Value = 8.5
$ws.Cells.Item(34, 23).Value = 8.75
$ws.Cells.Item(34, 25).Value = 28
$ws.Cells.Item(34, 27).Value = 9
$ws.Cells.Item(34, 28).Value = 20
$ws.Cells.Item(34, 29).Value = 100
$ws.Cells.Item(34, 30).Value = 18
$ws.Cells.Item(34, 31).Value = 45
$ws.Cells.Item(34, 32).Value = 22
$ws.Cells.Item(34, 34).Value = 75
$ws.Cells.Item(34, 36).Value = 800

# Row 35
$ws.Cells.Item(35, 8).Value = 3.55
$ws.Cells.Item(35, 20).Value = 6.8
$ws.Cells.Item(35, 21).Value = 7.6
$ws.Cells.Item(35, 24).Value = 13
$ws.Cells.Item(35, 28).Value = 15.5
$ws.Cells.Item(35, 30).Value = 14
$ws.Cells.Item(35, 34).Value = 50

# Row 36
$ws.Cells.Item(36, 12).Value = 1.32
$ws.Cells.Item(36, 20).Value = 7.6
$ws.Cells.Item(36, 25).Value = 27
$ws.Cells.Item(36, 31).Value = 16.5

# Row 37
$ws.Cells.Item(37, 7).Value = 1.72
$ws.Cells.Item(37, 8).Value = 3.4
$ws.Cells.Item(37, 9).Value = 4.4
$ws.Cells.Item(37, 21).Value = 7.6
$ws.Cells.Item(37, 23).Value = 13.5
$ws.Cells.Item(37, 27).Value = 6.8
$ws.Cells.Item(37, 28).Value = 17
$ws.Cells.Item(37, 33).Value = 75
$ws.Cells.Item(37, 34).Value = 45

# Row 38
$ws.Cells.Item(38, 8).Value = 3.35
$ws.Cells.Item(38, 24).Value = 18.5
$ws.Cells.Item(38, 34).Value = 24

# Row 42
$ws.Cells.Item(42, 7).Value = 2.12
$ws.Cells.Item(42, 8).Value = 3.3
$ws.Cells.Item(42, 9).Value = 3.15
$ws.Cells.Item(42, 12).Value = 1.37
$ws.Cells.Item(42, 13).Value = 2.65
$ws.Cells.Item(42, 14).Value = 2.07
$ws.Cells.Item(42, 15).Value = 1.6
$ws.Cells.Item(42, 17).Value = 2.45
$ws.Cells.Item(42, 18).Value = 1.88
$ws.Cells.Item(42, 19).Value = 1.72
$ws.Cells.Item(42, 20).Value = 6.6
$ws.Cells.Item(42, 21).Value = 9.5
$ws.Cells.Item(42, 22).Value = 9.25
$ws.Cells.Item(42, 23).Value = 19.5
$ws.Cells.Item(42, 24).Value = 19
$ws.Cells.Item(42, 25).Value = 35
$ws.Cells.Item(42, 27).Value = 6.4
$ws.Cells.Item(42, 30).Value = 8.5
$ws.Cells.Item(42, 31).Value = 15
$ws.Cells.Item(42, 32).Value = 11.5
$ws.Cells.Item(42, 33).Value = 40
$ws.Cells.Item(42, 34).Value = 30
$ws.Cells.Item(42, 35).Value = 45

